$d = $word.ActiveDocument

# 1. Locate the paragraph with the old "To be fixed" text and replace its
#    text directly via Range.Text (not Find&Replace) so that the apostrophe
#    in the new text is not mangled by smart-quote autocorrection.
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "To be fixed: ChooseLoginType screen has the wrong status bar color`r") {
        $targetParaIndex = $i
    }
}

$p = $d.Paragraphs.Item($targetParaIndex)
$pr = $p.Range
$textRange = $d.Range($pr.Start, $pr.End - 1)
$textRange.Text = "Solved: AnnotatedRegion isn't supposed to be at the top, Scaffold is the first one then "

# 2. Remove the now-redundant empty "NoSpacing" paragraph that immediately
#    follows the edited paragraph (the two empty paragraphs collapse to one).
$pAfter = $d.Paragraphs.Item($targetParaIndex + 1)
$pAfter.Range.Delete()

# 3. Remove the existing "_GoBack" bookmark that used to sit at the end of
#    the document (after the "On actual device keyboard..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 4. Re-create the "_GoBack" bookmark, collapsed, right after the run we
#    just edited (before its paragraph mark). A collapsed range passed
#    straight into Bookmarks.Add is not placed correctly, so insert a
#    placeholder character, bookmark the 1-character range, then delete
#    the placeholder -- this leaves a properly collapsed bookmark behind.
$p = $d.Paragraphs.Item($targetParaIndex)
$pr = $p.Range
$insertPos = $pr.End - 1
$insRange = $d.Range($insertPos, $insertPos)
$insRange.InsertAfter("X")
$charRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $charRange)
$delRange = $d.Range($insertPos, $insertPos + 1)
$delRange.Delete()
